$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for column C: "requiredMsg" ---
$ws.Cells.Item(1, 3).Value2 = "requiredMsg"
$ws.Cells.Item(1, 3).Font.Bold = $true
$ws.Cells.Item(1, 3).HorizontalAlignment = -4108
$ws.Cells.Item(1, 3).VerticalAlignment = -4108
$ws.Cells.Item(1, 3).WrapText = $true

# --- Data rows 2-9, column C: "Required" ---
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value2 = "Required"
    $ws.Cells.Item($r, 3).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 3).VerticalAlignment = -4108
}

# --- Data rows 2-9, column D: "Invalid credentials" ---
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "Invalid credentials"
    $ws.Cells.Item($r, 4).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 4).VerticalAlignment = -4108
}

# --- Header for column D: "invalidMsg" (added last) ---
$ws.Cells.Item(1, 4).Value2 = "invalidMsg"
$ws.Cells.Item(1, 4).Font.Bold = $true
$ws.Cells.Item(1, 4).HorizontalAlignment = -4108
$ws.Cells.Item(1, 4).VerticalAlignment = -4108
$ws.Cells.Item(1, 4).WrapText = $true

# --- Column widths for the new columns C and D ---
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws.Columns.Item(4).ColumnWidth = 17.666666666666668

# --- Update selection to match the new active cell/range ---
$ws.Range("D2:D9").Select() | Out-Null
